$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param($ws, $row, $col, $value)
    $ws.Cells.Item($row, $col).Value = $value
}

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 9
Set-CellValue $ws 9 8 549.5  # H9: 312 -> 549.5
Set-CellValue $ws 9 9 1000  # I9: 549.5 -> 1000
Set-CellValue $ws 9 10 99  # J9: 74.5 -> 99
Set-CellValue $ws 9 11 1000  # K9: 549.5 -> 1000
Set-CellValue $ws 9 12 99  # L9: 74.5 -> 99
Set-CellValue $ws 9 13 -831  # M9: -380.5 -> -831
Set-CellValue $ws 9 14 -437  # N9: -412.5 -> -437
# Row 15
Set-CellValue $ws 15 8 2861.5806  # H15: 3613.169 -> 2861.5806
Set-CellValue $ws 15 9 2861.5806  # I15: 3613.169 -> 2861.5806
Set-CellValue $ws 15 11 8584.7418  # K15: 10839.507 -> 8584.7418
Set-CellValue $ws 15 13 -8415.7418  # M15: -10670.507 -> -8415.7418
# Row 92
Set-CellValue $ws 92 8 808.9091  # H92: 879.8946999999999 -> 808.9091
Set-CellValue $ws 92 9 888.94116  # I92: 929.5625 -> 888.94116
Set-CellValue $ws 92 10 536.8  # J92: 615 -> 536.8
Set-CellValue $ws 92 11 888.94116  # K92: 929.5625 -> 888.94116
Set-CellValue $ws 92 12 536.8  # L92: 615 -> 536.8
Set-CellValue $ws 92 13 359.05884  # M92: 318.4375 -> 359.05884
Set-CellValue $ws 92 14 -3032.8  # N92: -3111 -> -3032.8
# Row 98
Set-CellValue $ws 98 8 4418.8184  # H98: 4441.8 -> 4418.8184
Set-CellValue $ws 98 9 4418.8184  # I98: 4441.8 -> 4418.8184
Set-CellValue $ws 98 11 4418.8184  # K98: 4441.8 -> 4418.8184
Set-CellValue $ws 98 13 -2920.8184  # M98: -2943.8 -> -2920.8184
# Row 122
Set-CellValue $ws 122 8 4418.8184  # H122: 4441.8 -> 4418.8184
Set-CellValue $ws 122 9 4418.8184  # I122: 4441.8 -> 4418.8184
Set-CellValue $ws 122 11 13256.4552  # K122: 13325.4 -> 13256.4552
Set-CellValue $ws 122 13 -10806.4552  # M122: -10875.4 -> -10806.4552
# Row 132
Set-CellValue $ws 132 8 9027.267  # H132: 9323.689 -> 9027.267
Set-CellValue $ws 132 9 6367.7085  # I132: 6625.826 -> 6367.7085
Set-CellValue $ws 132 11 19103.1255  # K132: 19877.478 -> 19103.1255
Set-CellValue $ws 132 13 -16573.1255  # M132: -17347.478 -> -16573.1255
# Row 137
Set-CellValue $ws 137 8 2005161.5  # H137: 2005169 -> 2005161.5
Set-CellValue $ws 137 9 2273592.5  # I137: 2381830.2 -> 2273592.5
Set-CellValue $ws 137 10 36668  # J137: 27698.25 -> 36668
Set-CellValue $ws 137 11 6820777.5  # K137: 7145490.600000001 -> 6820777.5
Set-CellValue $ws 137 12 110004  # L137: 83094.75 -> 110004
Set-CellValue $ws 137 13 -6818227.5  # M137: -7142940.600000001 -> -6818227.5
Set-CellValue $ws 137 14 -115104  # N137: -88194.75 -> -115104
# Row 141
Set-CellValue $ws 141 8 4012.8125  # H141: 4125.3125 -> 4012.8125
Set-CellValue $ws 141 9 2100.3572  # I141: 2169.6155 -> 2100.3572
Set-CellValue $ws 141 10 17400  # J141: 12600 -> 17400
Set-CellValue $ws 141 11 6301.071599999999  # K141: 6508.8465 -> 6301.071599999999
Set-CellValue $ws 141 12 52200  # L141: 37800 -> 52200
Set-CellValue $ws 141 13 -1121.071599999999  # M141: -1328.8465 -> -1121.071599999999
Set-CellValue $ws 141 14 -62560  # N141: -48160 -> -62560

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
Set-CellValue $ws 61 8 2527.6843  # H61: 2477.5642 -> 2527.6843
Set-CellValue $ws 61 9 1923.6154  # I61: 1854.6786 -> 1923.6154
Set-CellValue $ws 61 10 3836.5  # J61: 4063.0908 -> 3836.5
Set-CellValue $ws 61 11 1923.6154  # K61: 1854.6786 -> 1923.6154
Set-CellValue $ws 61 12 3836.5  # L61: 4063.0908 -> 3836.5
Set-CellValue $ws 61 13 -1711.6154  # M61: -1642.6786 -> -1711.6154
Set-CellValue $ws 61 14 -4260.5  # N61: -4487.0908 -> -4260.5
# Row 74
Set-CellValue $ws 74 8 254875.45  # H74: 233606.67 -> 254875.45
Set-CellValue $ws 74 9 398280.72  # I74: 328190 -> 398280.72
Set-CellValue $ws 74 10 3916.25  # J74: 3904.2856 -> 3916.25
Set-CellValue $ws 74 11 398280.72  # K74: 328190 -> 398280.72
Set-CellValue $ws 74 12 3916.25  # L74: 3904.2856 -> 3916.25
Set-CellValue $ws 74 13 -397406.72  # M74: -327316 -> -397406.72
Set-CellValue $ws 74 14 -5664.25  # N74: -5652.2856 -> -5664.25
# Row 77
Set-CellValue $ws 77 8 254875.45  # H77: 233606.67 -> 254875.45
Set-CellValue $ws 77 9 398280.72  # I77: 328190 -> 398280.72
Set-CellValue $ws 77 10 3916.25  # J77: 3904.2856 -> 3916.25
Set-CellValue $ws 77 11 1991403.6  # K77: 1640950 -> 1991403.6
Set-CellValue $ws 77 12 19581.25  # L77: 19521.428 -> 19581.25
Set-CellValue $ws 77 13 -1987035.6  # M77: -1636582 -> -1987035.6
Set-CellValue $ws 77 14 -28317.25  # N77: -28257.428 -> -28317.25
# Row 97
Set-CellValue $ws 97 8 2858.3333  # H97: 575 -> 2858.3333
Set-CellValue $ws 97 10 4000  # J97: 0 -> 4000
Set-CellValue $ws 97 12 4000  # L97: 0 -> 4000
Set-CellValue $ws 97 14 -4992  # N97: None -> -4992
# Row 102
Set-CellValue $ws 102 8 3174.4614  # H102: 2562.3572 -> 3174.4614
Set-CellValue $ws 102 9 2736.9  # I102: 1997.6364 -> 2736.9
Set-CellValue $ws 102 11 2736.9  # K102: 1997.6364 -> 2736.9
Set-CellValue $ws 102 13 -1114.9  # M102: -375.6364000000001 -> -1114.9
# Row 132
Set-CellValue $ws 132 8 2071.8684  # H132: 2158.861 -> 2071.8684
Set-CellValue $ws 132 9 1176.037  # I132: 1229.64 -> 1176.037
Set-CellValue $ws 132 11 3528.111  # K132: 3688.92 -> 3528.111
Set-CellValue $ws 132 13 -998.1109999999999  # M132: -1158.92 -> -998.1109999999999
# Row 136
Set-CellValue $ws 136 8 2527.6843  # H136: 2477.5642 -> 2527.6843
Set-CellValue $ws 136 9 1923.6154  # I136: 1854.6786 -> 1923.6154
Set-CellValue $ws 136 10 3836.5  # J136: 4063.0908 -> 3836.5
Set-CellValue $ws 136 11 5770.8462  # K136: 5564.0358 -> 5770.8462
Set-CellValue $ws 136 12 11509.5  # L136: 12189.2724 -> 11509.5
Set-CellValue $ws 136 13 -3220.8462  # M136: -3014.0358 -> -3220.8462
Set-CellValue $ws 136 14 -16609.5  # N136: -17289.2724 -> -16609.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
Set-CellValue $ws 86 8 2623.353  # H86: 2364.0476 -> 2623.353
Set-CellValue $ws 86 9 2429.2856  # I86: 2238.0625 -> 2429.2856
Set-CellValue $ws 86 10 3529  # J86: 2767.2 -> 3529
Set-CellValue $ws 86 11 2429.2856  # K86: 2238.0625 -> 2429.2856
Set-CellValue $ws 86 12 3529  # L86: 2767.2 -> 3529
Set-CellValue $ws 86 13 -1306.2856  # M86: -1115.0625 -> -1306.2856
Set-CellValue $ws 86 14 -5775  # N86: -5013.2 -> -5775
# Row 89
Set-CellValue $ws 89 8 2623.353  # H89: 2364.0476 -> 2623.353
Set-CellValue $ws 89 9 2429.2856  # I89: 2238.0625 -> 2429.2856
Set-CellValue $ws 89 10 3529  # J89: 2767.2 -> 3529
Set-CellValue $ws 89 11 12146.428  # K89: 11190.3125 -> 12146.428
Set-CellValue $ws 89 12 17645  # L89: 13836 -> 17645
Set-CellValue $ws 89 13 -6530.428  # M89: -5574.3125 -> -6530.428
Set-CellValue $ws 89 14 -28877  # N89: -25068 -> -28877
# Row 94
Set-CellValue $ws 94 8 250002600  # H94: 285716830 -> 250002600
Set-CellValue $ws 94 9 400001340  # I94: 500000960 -> 400001340
Set-CellValue $ws 94 11 400001340  # K94: 500000960 -> 400001340
Set-CellValue $ws 94 13 -400000889  # M94: -500000509 -> -400000889
# Row 134
Set-CellValue $ws 134 8 3385.5605  # H134: 3567.254 -> 3385.5605
Set-CellValue $ws 134 9 3163.0544  # I134: 3372.7925 -> 3163.0544
Set-CellValue $ws 134 10 4498.091  # J134: 4597.9 -> 4498.091
Set-CellValue $ws 134 11 9489.163199999999  # K134: 10118.3775 -> 9489.163199999999
Set-CellValue $ws 134 12 13494.273  # L134: 13793.7 -> 13494.273
Set-CellValue $ws 134 13 -6954.163199999999  # M134: -7583.377500000001 -> -6954.163199999999
Set-CellValue $ws 134 14 -18564.273  # N134: -18863.7 -> -18564.273

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
Set-CellValue $ws 16 8 1537.3793  # H16: 1552.1666 -> 1537.3793
Set-CellValue $ws 16 9 1376.0435  # I16: 1459.6666 -> 1376.0435
Set-CellValue $ws 16 10 2155.8333  # J16: 2199.6667 -> 2155.8333
Set-CellValue $ws 16 11 1376.0435  # K16: 1459.6666 -> 1376.0435
Set-CellValue $ws 16 12 2155.8333  # L16: 2199.6667 -> 2155.8333
Set-CellValue $ws 16 13 -1089.0435  # M16: -1172.6666 -> -1089.0435
Set-CellValue $ws 16 14 -2729.8333  # N16: -2773.6667 -> -2729.8333
# Row 31
Set-CellValue $ws 31 8 3390.7307  # H31: 3264.7036 -> 3390.7307
Set-CellValue $ws 31 9 2091.35  # I31: 2026.279 -> 2091.35
Set-CellValue $ws 31 10 7722  # J31: 8105.8184 -> 7722
Set-CellValue $ws 31 11 2091.35  # K31: 2026.279 -> 2091.35
Set-CellValue $ws 31 12 7722  # L31: 8105.8184 -> 7722
Set-CellValue $ws 31 13 -1796.35  # M31: -1731.279 -> -1796.35
Set-CellValue $ws 31 14 -8312  # N31: -8695.8184 -> -8312
# Row 34
Set-CellValue $ws 34 8 3390.7307  # H34: 3264.7036 -> 3390.7307
Set-CellValue $ws 34 9 2091.35  # I34: 2026.279 -> 2091.35
Set-CellValue $ws 34 10 7722  # J34: 8105.8184 -> 7722
Set-CellValue $ws 34 11 2091.35  # K34: 2026.279 -> 2091.35
Set-CellValue $ws 34 12 7722  # L34: 8105.8184 -> 7722
Set-CellValue $ws 34 13 -1889.35  # M34: -1824.279 -> -1889.35
Set-CellValue $ws 34 14 -8126  # N34: -8509.8184 -> -8126
# Row 58
Set-CellValue $ws 58 8 2928.88  # H58: 2813.037 -> 2928.88
Set-CellValue $ws 58 9 2206.2307  # I58: 2094.0667 -> 2206.2307
Set-CellValue $ws 58 11 2206.2307  # K58: 2094.0667 -> 2206.2307
Set-CellValue $ws 58 13 -2003.2307  # M58: -1891.0667 -> -2003.2307
# Row 113
Set-CellValue $ws 113 8 1537.3793  # H113: 1552.1666 -> 1537.3793
Set-CellValue $ws 113 9 1376.0435  # I113: 1459.6666 -> 1376.0435
Set-CellValue $ws 113 10 2155.8333  # J113: 2199.6667 -> 2155.8333
Set-CellValue $ws 113 11 1376.0435  # K113: 1459.6666 -> 1376.0435
Set-CellValue $ws 113 12 2155.8333  # L113: 2199.6667 -> 2155.8333
Set-CellValue $ws 113 13 793.9565  # M113: 710.3334 -> 793.9565
Set-CellValue $ws 113 14 -6495.8333  # N113: -6539.6667 -> -6495.8333
# Row 132
Set-CellValue $ws 132 8 2916.7932  # H132: 2869.5356 -> 2916.7932
Set-CellValue $ws 132 9 2482.4707  # I132: 2253.875 -> 2482.4707
Set-CellValue $ws 132 10 3532.0833  # J132: 3690.4167 -> 3532.0833
Set-CellValue $ws 132 11 7447.4121  # K132: 6761.625 -> 7447.4121
Set-CellValue $ws 132 12 10596.2499  # L132: 11071.2501 -> 10596.2499
Set-CellValue $ws 132 13 -4917.4121  # M132: -4231.625 -> -4917.4121
Set-CellValue $ws 132 14 -15656.2499  # N132: -16131.2501 -> -15656.2499
# Row 134
Set-CellValue $ws 134 8 2195.5757  # H134: 2371.5518 -> 2195.5757
Set-CellValue $ws 134 9 2165.9614  # I134: 2392.5454 -> 2165.9614
Set-CellValue $ws 134 11 6497.8842  # K134: 7177.6362 -> 6497.8842
Set-CellValue $ws 134 13 -3962.8842  # M134: -4642.6362 -> -3962.8842
# Row 136
Set-CellValue $ws 136 8 2928.88  # H136: 2813.037 -> 2928.88
Set-CellValue $ws 136 9 2206.2307  # I136: 2094.0667 -> 2206.2307
Set-CellValue $ws 136 11 6618.6921  # K136: 6282.2001 -> 6618.6921
Set-CellValue $ws 136 13 -4068.6921  # M136: -3732.2001 -> -4068.6921

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 87
Set-CellValue $ws 87 8 5990  # H87: 5990.143 -> 5990
Set-CellValue $ws 87 9 5990  # I87: 5990.143 -> 5990
Set-CellValue $ws 87 11 17970  # K87: 17970.429 -> 17970
Set-CellValue $ws 87 13 -16722  # M87: -16722.429 -> -16722
# Row 90
Set-CellValue $ws 90 8 5990  # H90: 5990.143 -> 5990
Set-CellValue $ws 90 9 5990  # I90: 5990.143 -> 5990
Set-CellValue $ws 90 11 53910  # K90: 53911.287 -> 53910
Set-CellValue $ws 90 13 -47670  # M90: -47671.287 -> -47670
# Row 123
Set-CellValue $ws 123 8 1855.138  # H123: 1851.7241 -> 1855.138
Set-CellValue $ws 123 9 999.9167  # I123: 991.6667 -> 999.9167
Set-CellValue $ws 123 11 2999.7501  # K123: 2975.0001 -> 2999.7501
Set-CellValue $ws 123 13 -549.7501000000002  # M123: -525.0001000000002 -> -549.7501000000002
# Row 130
Set-CellValue $ws 130 8 12530  # H130: 9140 -> 12530
Set-CellValue $ws 130 10 0  # J130: 5750 -> 0
Set-CellValue $ws 130 12 0  # L130: 17250 -> 0
$ws.Cells.Item(130, 14).ClearContents()  # N130 removed (was -27290)
# Row 136
Set-CellValue $ws 136 8 2303.1667  # H136: 2757.8 -> 2303.1667
Set-CellValue $ws 136 9 763.8  # I136: 947.25 -> 763.8
Set-CellValue $ws 136 11 2291.4  # K136: 2841.75 -> 2291.4
Set-CellValue $ws 136 13 2808.6  # M136: 2258.25 -> 2808.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
Set-CellValue $ws 97 8 2466.2307  # H97: 2638.4167 -> 2466.2307
Set-CellValue $ws 97 9 2253.111  # I97: 2484.75 -> 2253.111
Set-CellValue $ws 97 11 2253.111  # K97: 2484.75 -> 2253.111
Set-CellValue $ws 97 13 -1757.111  # M97: -1988.75 -> -1757.111
# Row 102
Set-CellValue $ws 102 8 2843.9092  # H102: 2999.375 -> 2843.9092
Set-CellValue $ws 102 9 2810.375  # I102: 2999.2 -> 2810.375
Set-CellValue $ws 102 10 2933.3333  # J102: 2999.6667 -> 2933.3333
Set-CellValue $ws 102 11 2810.375  # K102: 2999.2 -> 2810.375
Set-CellValue $ws 102 12 2933.3333  # L102: 2999.6667 -> 2933.3333
Set-CellValue $ws 102 13 -1188.375  # M102: -1377.2 -> -1188.375
Set-CellValue $ws 102 14 -6177.3333  # N102: -6243.6667 -> -6177.3333
# Row 132
Set-CellValue $ws 132 8 3248.2104  # H132: 3540.75 -> 3248.2104
Set-CellValue $ws 132 9 2576  # I132: 2780.923 -> 2576
Set-CellValue $ws 132 11 7728  # K132: 8342.769 -> 7728
Set-CellValue $ws 132 13 -5198  # M132: -5812.769 -> -5198
# Row 133
Set-CellValue $ws 133 8 208666.67  # H133: 176500 -> 208666.67
Set-CellValue $ws 133 10 208666.67  # J133: 176500 -> 208666.67
Set-CellValue $ws 133 12 208666.67  # L133: 176500 -> 208666.67
Set-CellValue $ws 133 14 -218786.67  # N133: -186620 -> -218786.67

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
Set-CellValue $ws 40 8 8046.222  # H40: 7697.923 -> 8046.222
Set-CellValue $ws 40 10 8884.5  # J40: 7899.375 -> 8884.5
Set-CellValue $ws 40 12 8884.5  # L40: 7899.375 -> 8884.5
Set-CellValue $ws 40 14 -9156.5  # N40: -8171.375 -> -9156.5
# Row 55
Set-CellValue $ws 55 8 245.5  # H55: 338.72726 -> 245.5
Set-CellValue $ws 55 9 67.85714  # I55: 170.85715 -> 67.85714
Set-CellValue $ws 55 10 660  # J55: 632.5 -> 660
Set-CellValue $ws 55 11 67.85714  # K55: 170.85715 -> 67.85714
Set-CellValue $ws 55 12 660  # L55: 632.5 -> 660
Set-CellValue $ws 55 13 105.14286  # M55: 2.14285000000001 -> 105.14286
Set-CellValue $ws 55 14 -1006  # N55: -978.5 -> -1006
# Row 103
Set-CellValue $ws 103 8 17999  # H103: 19999 -> 17999
Set-CellValue $ws 103 10 17999  # J103: 19999 -> 17999
Set-CellValue $ws 103 12 17999  # L103: 19999 -> 17999
Set-CellValue $ws 103 14 -20343  # N103: -22343 -> -20343
# Row 132
Set-CellValue $ws 132 8 7284.0557  # H132: 7813.3125 -> 7284.0557
Set-CellValue $ws 132 9 3836.5  # I132: 4229.75 -> 3836.5
Set-CellValue $ws 132 11 11509.5  # K132: 12689.25 -> 11509.5
Set-CellValue $ws 132 13 -8979.5  # M132: -10159.25 -> -8979.5
# Row 136
Set-CellValue $ws 136 8 6161.875  # H136: 5488.7 -> 6161.875
Set-CellValue $ws 136 9 8444  # I136: 6126.857 -> 8444
Set-CellValue $ws 136 10 3879.75  # J136: 3999.6667 -> 3879.75
Set-CellValue $ws 136 11 25332  # K136: 18380.571 -> 25332
Set-CellValue $ws 136 12 11639.25  # L136: 11999.0001 -> 11639.25
Set-CellValue $ws 136 13 -22782  # M136: -15830.571 -> -22782
Set-CellValue $ws 136 14 -16739.25  # N136: -17099.0001 -> -16739.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113
Set-CellValue $ws 113 8 462.33334  # H113: 543.2727 -> 462.33334
Set-CellValue $ws 113 9 357  # I113: 314.33334 -> 357
Set-CellValue $ws 113 10 989  # J113: 818 -> 989
Set-CellValue $ws 113 11 1071  # K113: 943.0000200000001 -> 1071
Set-CellValue $ws 113 12 2967  # L113: 2454 -> 2967
Set-CellValue $ws 113 13 1099  # M113: 1226.99998 -> 1099
Set-CellValue $ws 113 14 -7307  # N113: -6794 -> -7307
# Row 122
Set-CellValue $ws 122 8 11906165  # H122: 11906166 -> 11906165
Set-CellValue $ws 122 9 1284.7333  # I122: 1286.4667 -> 1284.7333
Set-CellValue $ws 122 11 3854.199900000001  # K122: 3859.4001 -> 3854.199900000001
Set-CellValue $ws 122 13 -1404.199900000001  # M122: -1409.4001 -> -1404.199900000001
# Row 132
Set-CellValue $ws 132 8 3069.3901  # H132: 3443.2 -> 3069.3901
Set-CellValue $ws 132 9 3145.3142  # I132: 3612.1724 -> 3145.3142
Set-CellValue $ws 132 11 9435.942599999998  # K132: 10836.5172 -> 9435.942599999998
Set-CellValue $ws 132 13 -6905.942599999998  # M132: -8306.5172 -> -6905.942599999998
# Row 140
Set-CellValue $ws 140 8 149000  # H140: 148974.5 -> 149000
Set-CellValue $ws 140 10 149000  # J140: 148974.5 -> 149000
Set-CellValue $ws 140 12 149000  # L140: 148974.5 -> 149000
Set-CellValue $ws 140 14 -159360  # N140: -159334.5 -> -159360
